$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.103.89'
$ws.Range('E2').Value = '  -3.77%  '
$ws.Range('D3').Value = '2.200.16'
$ws.Range('E3').Value = '  -3.59%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = "'105.69"
$ws.Range('E5').Value = '  -15.03%  '
$ws.Range('D6').Value = "'297.23"
$ws.Range('E6').Value = '  +11.35%  '
$ws.Range('D7').Value = "'0.619"
$ws.Range('E7').Value = '  -3.41%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = "'0.588"
$ws.Range('D10').Value = "'43.16"
$ws.Range('E10').Value = '  -10.85%  '
$ws.Range('D11').Value = "'0.0897"
$ws.Range('E11').Value = '  -5.76%  '
$ws.Range('D12').Value = "'53.98"
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = "'8.64"
$ws.Range('E13').Value = '  -7.88%  '
$ws.Range('D15').Value = "'0.925"
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = "'14.78"
$ws.Range('E16').Value = '  -4.58%  '
$ws.Range('D17').Value = '2.531.96'
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('D18').Value = '2.222.86'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').Value = '41.891.04'
$ws.Range('E19').Value = '  -4.15%  '
$ws.Range('D20').Value = "'7.28"
$ws.Range('E20').Value = '  +4.11%  '
$ws.Range('D21').Value = "'0.0000104"
$ws.Range('E21').Value = '  -6.15%  '
$ws.Range('D22').Value = "'71.98"
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').Value = "'3.51"
$ws.Range('E23').Value = '  +21.12%  '
$ws.Range('E24').Value = '  -7.91%  '
$ws.Range('D25').Value = "'226.04"
$ws.Range('E25').Value = '  -4.09%  '
$ws.Range('D26').Value = "'8.81"
$ws.Range('E26').Value = '  -7.55%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').Value = "'11.37"
$ws.Range('E28').Value = '  -4.75%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').Value = "'37.57"
$ws.Range('E31').Value = '  -11.06%  '
$ws.Range('E32').Value = '  -5.27%  '
$ws.Range('D33').Value = "'172.34"
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('D34').Value = "'20.68"
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('D35').Value = "'0.0866"
$ws.Range('E35').Value = '  -6.42%  '
$ws.Range('E36').Value = '  -5.61%  '
$ws.Range('D37').Value = "'4.90"
$ws.Range('E37').Value = '  +5.57%  '
$ws.Range('D38').Value = "'4.24"
$ws.Range('E38').Value = '  -2.58%  '
$ws.Range('E39').Value = '  -4.28%  '
$ws.Range('E40').Value = '  -5.44%  '
$ws.Range('E41').Value = '  -5.90%  '
$ws.Range('E42').Value = '  -5.12%  '
$ws.Range('D43').Value = "'69.79"
$ws.Range('E43').Value = '  -5.77%  '
$ws.Range('D44').Value = "'0.226"
$ws.Range('E44').Value = '  -5.41%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = "'12.45"
$ws.Range('E46').Value = '  -10.76%  '
$ws.Range('E47').Value = '  -6.91%  '
$ws.Range('E48').Value = '  -5.47%  '
$ws.Range('D49').Value = "'1.31"
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').Value = "'101.36"
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('E51').Value = '  -2.70%  '
